$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.303.15'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +7.60%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.018.86'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.25%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.76'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.16%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.79'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +7.20%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.017.43'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.23%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.518'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +2.46%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.40%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.450'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +4.27%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000242'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.41%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.27'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +6.91%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.90%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.128.82'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +7.34%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.509.95'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.96%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.97'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.20%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.010.18'
$ws.Range("D19").Style = "Normal"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '458.81'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.95%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.77'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +4.99%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.685'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.18%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.36'
$ws.Range("D23").Style = "Normal"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +7.68%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.77'
$ws.Range("D24").Style = "Normal"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.03%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.63'
$ws.Range("D25").Style = "Normal"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +4.50%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.26'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +11.76%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.76'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +7.44%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.42'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +17.67%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.87'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +11.71%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.61'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.22%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.29%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").Value = 'Hedera'

$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.112'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.75%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").Value = 'EthereumClassic'

$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.10'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.17%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.991'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.14%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.80'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +7.23%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.17'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +11.93%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +15.69%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.32'
$ws.Range("D40").Style = "Normal"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.96'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.50%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.26%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.302'
$ws.Range("D43").Style = "Normal"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +13.02%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.61%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '388.96'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +12.12%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.776.21'
$ws.Range("D46").Style = "Normal"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0353'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.33%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '135.68'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '23.48'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +8.55%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.14%  '
$ws.Range("E51").Style = "Normal"

